$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 2749
$ws1.Range("F11").Value = 348
$ws1.Range("F12").Value = 287
$ws1.Range("F14").Value = 5697
$ws1.Range("F15").Value = 618
$ws1.Range("F19").Value = 84
$ws1.Range("F20").Value = 459
$ws1.Range("F23").Value = 7

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 58
$ws2.Range("F13").Value = 633
$ws2.Range("F15").Value = 15

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 83
$ws3.Range("F5").Value = 2563
$ws3.Range("F9").Value = 1441

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2563
$ws4.Range("F7").Value = 1441
$ws4.Range("F15").Value = 2749
$ws4.Range("F21").Value = 348
$ws4.Range("F23").Value = 58
$ws4.Range("F25").Value = 5697
$ws4.Range("F28").Value = 618
$ws4.Range("F30").Value = 633
$ws4.Range("F33").Value = 84
$ws4.Range("F34").Value = 15
